# Update "想去人数" (want-to-go count) figures across the workbook's sheets.
# This mirrors a data refresh of the scraped 广州-漫展信息 workbook
# (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9
$ws1.Range("F4").Value = 1325
$ws1.Range("F7").Value = 3921
$ws1.Range("F9").Value = 783
$ws1.Range("F10").Value = 2352
$ws1.Range("F12").Value = 48
$ws1.Range("F14").Value = 759
$ws1.Range("F16").Value = 194
$ws1.Range("F17").Value = 2529
$ws1.Range("F19").Value = 234
$ws1.Range("F22").Value = 237

# --- 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 56
$ws2.Range("F4").Value = 29
$ws2.Range("F18").Value = 53

# --- 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 2123
$ws3.Range("F5").Value = 346
$ws3.Range("F6").Value = 16

# --- 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2123
$ws4.Range("F5").Value = 346
$ws4.Range("F6").Value = 56
$ws4.Range("F7").Value = 56
$ws4.Range("F8").Value = 9
$ws4.Range("F11").Value = 1325
$ws4.Range("F13").Value = 29
$ws4.Range("F16").Value = 16
$ws4.Range("F18").Value = 3921
$ws4.Range("F24").Value = 783
$ws4.Range("F25").Value = 2352
$ws4.Range("F27").Value = 48
$ws4.Range("F30").Value = 759
$ws4.Range("F32").Value = 194
$ws4.Range("F38").Value = 234
$ws4.Range("F41").Value = 237
$ws4.Range("F45").Value = 53
